$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns C:D (128 Channels, 45 deg FOV and 128 Channels, 22.5 deg FOV)
$ws.Range("C1:D9").Delete()

# Update the "angle offset" value for row 2 (64 Channels, 45 deg FOV) from 6 to 0
$ws.Range("B2").Value = 0
